{"js": "// Apply benchmark-table updates for the Shenandoah GC lusearch heap-1G docs.\n// The document is a single-column table; each row holds one benchmark\n// figure (a handful of rows hold a whole tab-separated line in one run).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index (0-based) -> new text for that cell\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"1000\"],\n  [4, \"0.00002\"],\n  [6, \"0.00015\"],\n  [8, \"0.00022\"],\n  [9, \"0.00038\"],\n  [10, \"0.00049\"],\n  [11, \"0.19567\"],\n  [43, \"99.88\"],\n  [44, \"0.2\"],\n  [45, \"162\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  const range = cell.body.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply benchmark-table updates for the Shenandoah GC lusearch heap-1G docs.\n# The document is a single-column table; each row holds one benchmark\n# figure (a handful of rows hold a whole tab-separated line in one run).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Word COM row/column indices are 1-based.\n$tbl.Cell(1, 1).Range.Text = \"0M\"\n$tbl.Cell(2, 1).Range.Text = \"0M\"\n$tbl.Cell(3, 1).Range.Text = \"0M\"\n$tbl.Cell(4, 1).Range.Text = \"1000\"\n$tbl.Cell(5, 1).Range.Text = \"0.00002\"\n$tbl.Cell(7, 1).Range.Text = \"0.00015\"\n$tbl.Cell(9, 1).Range.Text = \"0.00022\"\n$tbl.Cell(10, 1).Range.Text = \"0.00038\"\n$tbl.Cell(11, 1).Range.Text = \"0.00049\"\n$tbl.Cell(12, 1).Range.Text = \"0.19567\"\n$tbl.Cell(44, 1).Range.Text = \"99.88\"\n$tbl.Cell(45, 1).Range.Text = \"0.2\"\n$tbl.Cell(46, 1).Range.Text = \"162\"\n"}
